# Apply the scraped price/volume refresh to the cryptos sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal text value, forcing text storage (leading
# apostrophe, exactly like a user typing '301.84 into Excel) so numeric-looking
# strings such as "301.84" are not silently reinterpreted as numbers.
function Set-TextValue($range, $text) {
    $ws.Range($range).Value = "'$text"
}

$ws.Range("D2").Value = "46.680.52"
$ws.Range("E2").Value = "  +6.19%  "
$ws.Range("D3").Value = "2.305.96"
$ws.Range("E3").Value = "  +4.94%  "
$ws.Range("E4").Value = "  -0.48%  "
Set-TextValue "D5" "301.84"
$ws.Range("E5").Value = "  +2.10%  "
Set-TextValue "D6" "101.99"
$ws.Range("E6").Value = "  +13.77%  "
Set-TextValue "D7" "0.574"
$ws.Range("E7").Value = "  +1.78%  "
Set-TextValue "D8" "0.999"
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("E9").Value = "  +8.90%  "
Set-TextValue "D10" "36.94"
$ws.Range("E10").Value = "  +13.90%  "
Set-TextValue "D11" "0.0806"
$ws.Range("E11").Value = "  +4.09%  "
Set-TextValue "D12" "7.40"
$ws.Range("E12").Value = "  +7.77%  "
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("D14").Value = "2.657.48"
$ws.Range("E14").Value = "  +5.00%  "
$ws.Range("D15").Value = "2.307.53"
$ws.Range("E15").Value = "  +2.09%  "
Set-TextValue "D16" "14.04"
$ws.Range("E16").Value = "  +5.82%  "
Set-TextValue "D17" "0.823"
$ws.Range("E17").Value = "  +5.77%  "
$ws.Range("D18").Value = "46.611.26"
$ws.Range("E18").Value = "  +6.70%  "
Set-TextValue "D19" "13.26"
$ws.Range("E19").Value = "  +21.04%  "
$ws.Range("E20").Value = "  +5.97%  "
Set-TextValue "D21" "6.15"
$ws.Range("E21").Value = "  +5.40%  "
Set-TextValue "D22" "66.87"
$ws.Range("E22").Value = "  +5.53%  "
Set-TextValue "D23" "248.84"
$ws.Range("E23").Value = "  +7.01%  "
$ws.Range("E24").Value = "  +6.40%  "
$ws.Range("E25").Value = "  +6.65%  "
$ws.Range("E26").Value = "  -1.08%  "
Set-TextValue "D27" "43.59"
$ws.Range("E27").Value = "  +18.77%  "
Set-TextValue "D28" "2.28"
$ws.Range("E28").Value = "  +1.58%  "
Set-TextValue "D29" "9.92"
$ws.Range("E29").Value = "  +6.82%  "
Set-TextValue "D30" "20.17"
$ws.Range("E30").Value = "  +4.43%  "
Set-TextValue "D31" "5.83"
$ws.Range("E31").Value = "  +9.40%  "
Set-TextValue "D32" "0.0805"
$ws.Range("E32").Value = "  +9.05%  "
Set-TextValue "D33" "146.48"
$ws.Range("E33").Value = "  -1.43%  "
$ws.Range("E34").Value = "  +4.23%  "
Set-TextValue "D35" "3.18"
$ws.Range("E35").Value = "  +11.91%  "
Set-TextValue "D36" "0.111"
$ws.Range("E36").Value = "  +8.03%  "
$ws.Range("E37").Value = "  +3.38%  "
$ws.Range("E38").Value = "  +10.58%  "
Set-TextValue "D39" "15.59"
$ws.Range("E39").Value = "  +18.84%  "
Set-TextValue "D40" "4.08"
$ws.Range("E40").Value = "  +15.54%  "
Set-TextValue "D41" "3.49"
$ws.Range("E41").Value = "  +12.93%  "
$ws.Range("E42").Value = "  +6.39%  "
Set-TextValue "D43" "0.998"
$ws.Range("E43").Value = "  -0.63%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.857.10"
$ws.Range("E44").Value = "  +3.11%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D45" "1.97"
$ws.Range("E45").Value = "  +17.02%  "
Set-TextValue "D46" "88.38"
$ws.Range("E46").Value = "  +21.70%  "
Set-TextValue "D47" "0.197"
$ws.Range("E47").Value = "  +12.32%  "
Set-TextValue "D48" "74.59"
$ws.Range("E48").Value = "  +14.26%  "
Set-TextValue "D49" "4.93"
$ws.Range("E49").Value = "  +11.98%  "
Set-TextValue "D50" "97.65"
$ws.Range("E50").Value = "  +5.89%  "
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D51" "8.08"
$ws.Range("E51").Value = "  +8.02%  "
